$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("E2").Value = "Andre B.-Eletrônica analóg. e de potên"

$ws.Range("B3").Value = "-"
$ws.Range("E3").Value = "[-, André Guimarães-CAD]"

$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "[-, Sandro-Programação de Computadores, Pedro Bispo-Acionamentos Elétricos, -]"
$ws.Range("E4").Value = "[-, André Guimarães-CAD]"
$ws.Range("F4").Value = "Sandro-Circuitos elétrico"

$ws.Range("B6").Value = "[João Paulo-Sistemas digitais, -, -]"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "[-, Sandro-Programação de Computadores, Pedro Bispo-Acionamentos Elétricos, -]"
$ws.Range("E6").Value = "Andre B.-Eletrônica analóg. e de potên"
$ws.Range("F6").Value = "Sandro-Circuitos elétrico"

$ws.Range("B7").Value = "[João Paulo-Sistemas digitais, -, -]"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "[-, Sandro-Programação de Computadores, Pedro Bispo-Acionamentos Elétricos, -]"
$ws.Range("E7").Value = "Nilton Maia-M.T.R"

$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "[-, Sandro-Programação de Computadores, Pedro Bispo-Acionamentos Elétricos, -]"
$ws.Range("E8").Value = "Nilton Maia-M.T.R"
